# Applies commit "updated statbar xpaths & diagnosis testcases":
# adds CypherOutput_Message, StatOutput, and StatOutput_Message sheets
# that mirror the existing Message sheet plus a new stats query/result block.

$wb = $excel.ActiveWorkbook

$neo4jUrlLabel   = "Neo4j_URL:"
$neo4jUrlValue   = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$userLabel       = "User_name:"
$userValue       = "neo4j"
$pwdLabel        = "PWD:"
$pwdValue        = "icdcDBneo4j0"
$cypherLabel     = "Cypher:"
$cypherQueryOld  = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Lip and oropharyngeal neoplasms malignant :: Melanoma-mucosa/mandible''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$cypherQueryNew  = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Lip and oropharyngeal neoplasms malignant :: Melanoma-mucosa/mandible'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$outputLabel     = "Output:"
$outputPath      = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC03_Canine_Filter_Diagnosis-LipMandible_Neo4jData.xlsx'

function Write-MessageBlock($ws, $startRow, $cypherQuery) {
    $ws.Cells.Item($startRow, 1).Value = $neo4jUrlLabel
    $ws.Cells.Item($startRow + 1, 1).Value = $neo4jUrlValue
    $ws.Cells.Item($startRow + 2, 1).Value = $userLabel
    $ws.Cells.Item($startRow + 3, 1).Value = $userValue
    $ws.Cells.Item($startRow + 4, 1).Value = $pwdLabel
    $ws.Cells.Item($startRow + 5, 1).Value = $pwdValue
    $ws.Cells.Item($startRow + 6, 1).Value = $cypherLabel
    $ws.Cells.Item($startRow + 7, 1).Value = $cypherQuery
    $ws.Cells.Item($startRow + 8, 1).Value = $outputLabel
    $ws.Cells.Item($startRow + 9, 1).Value = $outputPath
}

$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)

# --- New sheet: CypherOutput_Message (duplicate of Message) ---
$sCypherMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sCypherMsg.Name = "CypherOutput_Message"
Write-MessageBlock $sCypherMsg 1 $cypherQueryOld

# --- New sheet: StatOutput (counts returned by the new stats cypher query) ---
$sStatOutput = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sCypherMsg)
$sStatOutput.Name = "StatOutput"
$sStatOutput.Cells.Item(1,1).Value = "number_of_files"
$sStatOutput.Cells.Item(1,2).Value = "number_of_sample"
$sStatOutput.Cells.Item(1,3).Value = "number_of_cases"
$sStatOutput.Cells.Item(1,4).Value = "number_of_study"
# Store the numeric results as text so they are written as shared strings,
# matching the CypherOutput/Message convention used throughout this workbook.
$sStatOutput.Range("A2:D2").NumberFormat = "@"
$sStatOutput.Cells.Item(2,1).Value = "187"
$sStatOutput.Cells.Item(2,2).Value = "17"
$sStatOutput.Cells.Item(2,3).Value = "7"
$sStatOutput.Cells.Item(2,4).Value = "1"

# --- New sheet: StatOutput_Message (Message block repeated for old + new query) ---
$sStatMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sStatOutput)
$sStatMsg.Name = "StatOutput_Message"
Write-MessageBlock $sStatMsg 1 $cypherQueryOld
Write-MessageBlock $sStatMsg 11 $cypherQueryNew

# Restore original active sheet/tab selection.
$wb.Worksheets.Item(1).Activate()
